$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: merge the two runs "SUN Dec 23" + " 11:27:09 IST 2018" into
# a single run "SUN Dec 23 11:27:09 IST 2018" (text itself is unchanged,
# only the run split disappears). A Find/Replace over the same text
# causes Word to collapse the matched range into one run.
# ---------------------------------------------------------------------
$rng1 = $d.Content
$rng1.Find.ClearFormatting()
$rng1.Find.Execute("SUN Dec 23 11:27:09 IST 2018", $false, $false, $false, $false, $false, `
                    $true, 1, $false, "SUN Dec 23 11:27:09 IST 2018", 2) | Out-Null

# ---------------------------------------------------------------------
# Change 2: append a brand-new purchase entry (THU Jan 10 12:03:12 IST
# 2019 / LGL B / BEET / 1620 / CASH AND CLEARD) after the last existing
# entry, right after the paragraph that holds the final "- CASH" text,
# followed by two blank paragraphs (mirroring the template used by the
# other entries in the document).
# ---------------------------------------------------------------------
$w_ns = "xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`""

$pStyleNormal = '<w:pPr><w:pStyle w:val="PlainText"/><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr></w:pPr>'
$pStyleRed    = '<w:pPr><w:pStyle w:val="PlainText"/><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="FF0000"/></w:rPr></w:pPr>'

$rPrNormal = '<w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr>'
$rPrRed    = '<w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="FF0000"/></w:rPr>'

# paragraph 1: blank line
$p1 = "<w:p $w_ns>$pStyleNormal</w:p>"

# paragraph 2: timestamp "THU Jan 10" / " 12:03:12 IST 2019"
$p2 = "<w:p $w_ns>$pStyleNormal" + `
        "<w:r>$rPrNormal<w:t>THU Jan 10</w:t></w:r>" + `
        "<w:r>$rPrNormal<w:t xml:space=`"preserve`"> 12:03:12 IST 2019</w:t></w:r>" + `
      "</w:p>"

# paragraph 3: "Person Name" ... "- LGL B"
$p3 = "<w:p $w_ns>$pStyleNormal" + `
        "<w:r>$rPrNormal<w:t>Person Name</w:t></w:r>" + `
        "<w:r>$rPrNormal<w:tab/></w:r>" + `
        "<w:r>$rPrNormal<w:tab/></w:r>" + `
        "<w:r>$rPrNormal<w:tab/></w:r>" + `
        "<w:r>$rPrNormal<w:tab/><w:t>- LGL B</w:t></w:r>" + `
      "</w:p>"

# paragraph 4: separator dashes
$p4 = "<w:p $w_ns>$pStyleNormal" + `
        "<w:r>$rPrNormal<w:t>---------------------------------------------------------------</w:t></w:r>" + `
      "</w:p>"

# paragraph 5: "Item Name" ... "- BEET"
$p5 = "<w:p $w_ns>$pStyleNormal" + `
        "<w:r>$rPrNormal<w:t>Item Name</w:t></w:r>" + `
        "<w:r>$rPrNormal<w:tab/></w:r>" + `
        "<w:r>$rPrNormal<w:tab/></w:r>" + `
        "<w:r>$rPrNormal<w:tab/></w:r>" + `
        "<w:r>$rPrNormal<w:tab/><w:t>- BEET</w:t></w:r>" + `
      "</w:p>"

# paragraph 6: "Amount Received" (red) ... "- 1620" (red)
$p6 = "<w:p $w_ns>$pStyleRed" + `
        "<w:r>$rPrRed<w:t>Amount Received</w:t></w:r>" + `
        "<w:r>$rPrRed<w:tab/></w:r>" + `
        "<w:r>$rPrRed<w:tab/></w:r>" + `
        "<w:r>$rPrRed<w:tab/><w:t>- 1620</w:t></w:r>" + `
      "</w:p>"

# paragraph 7: "Amount Received mode" ... "- CASH AND CLEARD"
$p7 = "<w:p $w_ns>$pStyleNormal" + `
        "<w:r>$rPrNormal<w:t>Amount Received mode</w:t></w:r>" + `
        "<w:r>$rPrNormal<w:tab/></w:r>" + `
        "<w:r>$rPrNormal<w:tab/><w:t>- CASH AND CLEARD</w:t></w:r>" + `
      "</w:p>"

# paragraph 8: blank line
$p8 = "<w:p $w_ns>$pStyleNormal</w:p>"

# paragraph 9: reconstruction of the blank paragraph that originally
# followed the final "- CASH" line (so it is preserved, just pushed down)
$p9 = "<w:p $w_ns>$pStyleNormal</w:p>"

$newXml = $p1 + $p2 + $p3 + $p4 + $p5 + $p6 + $p7 + $p8 + $p9

# Locate the paragraph that used to directly follow the last "- CASH"
# paragraph (i.e. the final blank paragraph of the document template)
# and replace its (collapsed) start with the new block followed by a
# reconstruction of itself - this has the net effect of inserting the
# new paragraphs right after the "- CASH" paragraph while leaving the
# pre-existing blank paragraph intact just after them.
$lastCash = $null
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $ptext = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13)
    if ($ptext.StartsWith("Amount Received mode") -and $ptext.EndsWith("- CASH")) {
        $lastCash = $i
        break
    }
}

$targetPara = $d.Paragraphs.Item($lastCash + 1)
$targetRange = $targetPara.Range
$targetRange.Collapse(1)
$targetRange.InsertXML($newXml)
